$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.248.78"
Set-TextValue "E2" "  -0.36%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.287.81"
Set-TextValue "E3" "  +0.32%  "

# Row 5 - BNB
Set-TextValue "D5" "535.48"
Set-TextValue "E5" "  -1.84%  "

# Row 6 - Solana
Set-TextValue "D6" "131.16"
Set-TextValue "E6" "  +0.29%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.02%  "

# Row 8 - XRP
Set-TextValue "E8" "  +3.03%  "

# Row 9 - LidoStakedEther
Set-TextValue "D9" "2.282.88"
Set-TextValue "E9" "  +0.16%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.1000"
Set-TextValue "E10" "  -1.39%  "

# Row 11 - Toncoin
Set-TextValue "D11" "5.47"
Set-TextValue "E11" "  -0.76%  "

# Row 12 - TRON
Set-TextValue "E12" "  +1.02%  "

# Row 13 - Cardano
Set-TextValue "E13" "  -0.32%  "

# Row 14 - Avalanche
Set-TextValue "D14" "23.57"
Set-TextValue "E14" "  -0.15%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "2.694.28"
Set-TextValue "E15" "  +0.30%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "58.165.64"

# Row 17 - ShibaInu
Set-TextValue "E17" "  -0.38%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "2.272.30"
Set-TextValue "E18" "  -0.38%  "

# Row 19 - Chainlink
Set-TextValue "E19" "  -0.75%  "

# Row 20 - Polkadot
Set-TextValue "E20" "  -2.71%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "313.28"
Set-TextValue "E21" "  -0.35%  "

# Row 22 - Uniswap
Set-TextValue "E22" "  +1.22%  "

# Row 23 - Dai
Set-TextValue "E23" "  +0.16%  "

# Row 24 - Litecoin
Set-TextValue "D24" "63.19"
Set-TextValue "E24" "  +0.31%  "

# Row 25 - Kaspa
Set-TextValue "D25" "0.167"
Set-TextValue "E25" "  -1.37%  "

# Row 26 - Binance-PegBSC-USD
Set-TextValue "E26" "  +0.09%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "7.99"
Set-TextValue "E27" "  -1.40%  "

# Row 28 - Fetch.AI
Set-TextValue "E28" "  -1.18%  "

# Row 29 - Monero
Set-TextValue "D29" "170.76"
Set-TextValue "E29" "  +0.10%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  -2.48%  "

# Row 31 - PEPE
Set-TextValue "E31" "  +0.18%  "

# Row 32 - Aptos
Set-TextValue "D32" "5.81"
Set-TextValue "E32" "  +1.17%  "

# Row 33 - SuiNetwork
Set-TextValue "D33" "1.07"
Set-TextValue "E33" "  -1.25%  "

# Row 34 - PolygonEcosystemToken
Set-TextValue "D34" "0.380"
Set-TextValue "E34" "  -0.88%  "

# Row 36 - EthereumClassic
Set-TextValue "D36" "17.83"
Set-TextValue "E36" "  +0.30%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "E37" "  -0.08%  "

# Row 38 - ImmutableX
Set-TextValue "E38" "  -1.18%  "

# Row 39 - NEARProtocol
Set-TextValue "E39" "  -0.71%  "

# Row 40 - Stacks
Set-TextValue "E40" "  -1.52%  "

# Row 41 - Bittensor
Set-TextValue "D41" "289.04"
Set-TextValue "E41" "  -3.59%  "

# Row 42 - Aave
Set-TextValue "D42" "140.44"
Set-TextValue "E42" "  -0.52%  "

# Row 43 - Filecoin
Set-TextValue "E43" "  -0.28%  "

# Row 44 - Stellar
Set-TextValue "D44" "0.0950"
Set-TextValue "E44" "  +0.31%  "

# Row 45 - Hedera
Set-TextValue "D45" "0.0494"
Set-TextValue "E45" "  -0.77%  "

# Row 46 - Mantle
Set-TextValue "D46" "0.553"
Set-TextValue "E46" "  +0.06%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "18.10"
Set-TextValue "E47" "  -1.50%  "

# Row 48 - VeChain
Set-TextValue "E48" "  -1.67%  "

# Row 49 - WhiteBITCoin
Set-TextValue "E49" "  -0.47%  "

# Row 50 - ZEEBU
Set-TextValue "E50" "  +0.27%  "

# Row 51 - dogwifhat
Set-TextValue "E51" "  +1.46%  "
